$d = $word.ActiveDocument

# --- Remove the stray _GoBack bookmark from the first (empty) paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Add a new paragraph after "Applicant Signature " with a line break,
#     "Dated: {{ today() }}" text, and re-create the _GoBack bookmark there ---
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newRange = $newPara.Range

$datedXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/part1.xml" pkg:contentType="application/xml"><pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p><w:r><w:br/><w:t>Dated: {{ today() }}</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$newRange.InsertXML($datedXml)

# InsertXML leaves one extra empty trailing paragraph behind (the original
# paragraph mark of the freshly-inserted empty paragraph) - remove it.
$trailing = $d.Paragraphs.Item($d.Paragraphs.Count)
$cleanupRange = $d.Range($trailing.Range.Start - 1, $trailing.Range.End)
$cleanupRange.Delete()
